$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 7: fill in the new time log entry
$ws.Range("A7").Value = Get-Date -Year 2014 -Month 10 -Day 21 -Hour 0 -Minute 0 -Second 0
$ws.Range("B7").Value = "9:00am"
$ws.Range("C7").Value = "11:20am"
$ws.Range("D7").Value = 140
$ws.Range("E7").Value = "Programming"
$ws.Range("F7").Value = "Worked on BreakOut Collision Ball and Player"
